$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SCENARIO_DESC (column D) text for rows 2-5 with shorter descriptions
$ws.Range("D2").Value = "Tambah Setup Sektor"
$ws.Range("D3").Value = "View Setup Sektor"
$ws.Range("D4").Value = "Ubah Setup Sektor"
$ws.Range("D5").Value = "Hapus Setup Sektor"

# Adjust row heights to fit the shorter text
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).RowHeight = 30

# Update the selected cell
$ws.Range("D5").Select()
